# Fix Training Data Issue (#48)
#
# The "Date" column (BF) stored values like "6-29-2007-08" - a mash-up of the
# game's month/day ("6-29") and the season label ("2007-08"). Because of how
# NBA.com reports stats, games played on June 29 during the "2007-08" season
# actually happened in calendar year 2008, so the correct ISO date is
# 2008-06-29. This script rewrites every row's Date cell accordingly.
#
# Note: assigning a literal ISO-looking string straight to Value/Value2 makes
# Excel's input parser interpret it as a date serial (and reformat the cell),
# which is not what we want here - we need the literal text "2008-06-29".
# Building the text via a formula and pasting just the *value* preserves it
# as plain text without touching any cell's number format/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "6-29-2007-08"
$newValue = "2008-06-29"

# Off-sheet scratch cell used purely to materialize the literal text value.
$helper = $ws.Cells.Item(1000, 1)
$helper.Formula = '="' + $newValue + '"'

$firstRow = 2
$lastRow = 31
$dateCol = 58  # column BF

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldValue) {
        $helper.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

# Clean up the scratch cell and the clipboard marching-ants state.
$helper.Value2 = ""
$excel.CutCopyMode = $false
